$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = [double]"0.9999999984649723"
$ws.Range("C2").Value = [double]"5.150631097231949e-10"
$ws.Range("D2").Value = [double]"6.599353474689541e-10"
$ws.Range("E2").Value = [double]"3.600186320736702e-10"
$ws.Range("F2").Value = 46059
